# Applies the "Further minor improvements and labelling" commit.
#
# Summary of intentional content changes (per the OOXML diff):
#  1. Sheet1 renamed "Timesatser_budget" -> "Timesatser"
#  2. Sheet1 becomes the active/selected tab (was sheet "Eksterne timer og
#     porteføljer"); new selection on Sheet1 is F27.
#  3. Sheet1 B10:F10 (the budget totals row) reformatted from 2-decimal to
#     a whole-number ("0") number format.
#  4. Sheet2 ("Preallokering") row 2 re-pointed: allocation moves from
#     columns B/C to columns B/E; new selection on Sheet2 is E3.
#  5. Sheet3 ("Eksterne timer og porteføljer") row 2, column B value
#     cleared to 0; new selection on Sheet3 is B3.
#
# (Sub-pixel column-width / default-row-height / dyDescent churn and the
# absPath / revisionPtr GUIDs in the diff are environment/render artifacts
# from the authoring machine - not addressable, or meaningful, via the
# Excel object model - and are intentionally left alone.)

$wb = $excel.ActiveWorkbook

# --- 1. Rename first sheet ---------------------------------------------
$wsRates = $wb.Worksheets.Item(1)
$wsRates.Name = "Timesatser"

$wsPre = $wb.Worksheets.Item("Preallokering")
$wsExt = $wb.Worksheets.Item("Eksterne timer og porteføljer")

# --- 3. Whole-number formatting for the budget totals row -------------
$wsRates.Range("B10:F10").NumberFormat = "0"

# --- 4. Preallokering row 2 data edits ----------------------------------
$wsPre.Range("B2").Value = 200
$wsPre.Range("C2").Value = 0
$wsPre.Range("D2").Value = 0
$wsPre.Range("E2").Value = 50
$wsPre.Range("F2").Value = 0

# --- 5. Eksterne timer og porteføljer row 2 data edit -------------------
$wsExt.Range("B2").Value = 0

# --- 2. Selections, then make Sheet1 the active tab last ---------------
$wsPre.Range("E3").Select()
$wsExt.Range("B3").Select()

$wsRates.Activate()
$wsRates.Range("F27").Select()
